$wb = $excel.ActiveWorkbook

# New sensor-log rows (rows 73-85) to append to the PIR, Humidity and
# Temperature sheets. Columns are: Date, Timestamp, Hour, Location, Value, Status

$pirRows = @(
    @("2026-01-28","17:07:58","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:07:59","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:04","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:09","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:14","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:19","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:24","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:29","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:34","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:39","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:44","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:49","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","17:08:54","17:00","Bathroom","No Motion","Inactive")
)

$humidityRows = @(
    @("2026-01-28","17:07:57","17:00","Bathroom","87.6%","Active"),
    @("2026-01-28","17:08:00","17:00","Bathroom","87.5%","Active"),
    @("2026-01-28","17:08:12","17:00","Bathroom","87.5%","Active"),
    @("2026-01-28","17:08:16","17:00","Bathroom","86.6%","Active"),
    @("2026-01-28","17:08:20","17:00","Bathroom","87.5%","Active"),
    @("2026-01-28","17:08:24","17:00","Bathroom","86.7%","Active"),
    @("2026-01-28","17:08:28","17:00","Bathroom","87.5%","Active"),
    @("2026-01-28","17:08:32","17:00","Bathroom","87.5%","Active"),
    @("2026-01-28","17:08:36","17:00","Bathroom","86.6%","Active"),
    @("2026-01-28","17:08:40","17:00","Bathroom","87.5%","Active"),
    @("2026-01-28","17:08:48","17:00","Bathroom","86.7%","Active"),
    @("2026-01-28","17:08:53","17:00","Bathroom","87.5%","Active"),
    @("2026-01-28","17:08:57","17:00","Bathroom","86.7%","Active")
)

$temperatureRows = @(
    @("2026-01-28","17:07:57","17:00","Bathroom","22.9C","Active"),
    @("2026-01-28","17:08:01","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:13","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:17","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:21","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:25","17:00","Bathroom","22.9C","Active"),
    @("2026-01-28","17:08:29","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:33","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:37","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:41","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:49","17:00","Bathroom","22.9C","Active"),
    @("2026-01-28","17:08:53","17:00","Bathroom","22.8C","Active"),
    @("2026-01-28","17:08:57","17:00","Bathroom","22.9C","Active")
)

function Write-SheetRows {
    param($sheetName, $rows, $startRow)

    $ws = $wb.Worksheets.Item($sheetName)
    $r = $startRow
    foreach ($row in $rows) {
        for ($c = 0; $c -lt $row.Length; $c++) {
            $value = $row[$c]
            # Force values that look like dates/percentages/numbers to stay
            # plain text, matching the source log's text-only columns.
            if ($value -match "^\d{4}-\d{2}-\d{2}$" -or $value.EndsWith("%")) {
                $ws.Cells.Item($r, $c + 1).Value = "'" + $value
            } else {
                $ws.Cells.Item($r, $c + 1).Value = $value
            }
        }
        $r = $r + 1
    }
}

Write-SheetRows "PIR" $pirRows 73
Write-SheetRows "Humidity" $humidityRows 73
Write-SheetRows "Temperature" $temperatureRows 73

Write-Host "Appended rows 73-85 to PIR, Humidity and Temperature sheets"
